$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 154
$ws.Range("F7").Value = 156
$ws.Range("F8").Value = 754
$ws.Range("F9").Value = 4125
$ws.Range("F12").Value = 166
$ws.Range("F14").Value = 5892
$ws.Range("F16").Value = 2284
$ws.Range("F18").Value = 157
$ws.Range("F20").Value = 8862
$ws.Range("F22").Value = 1695
$ws.Range("F23").Value = 193
$ws.Range("F24").Value = 2277
$ws.Range("F25").Value = 2371
$ws.Range("F26").Value = 1373
$ws.Range("F28").Value = 1922
$ws.Range("F30").Value = 52
$ws.Range("F33").Value = 34
$ws.Range("F35").Value = 37
$ws.Range("F36").Value = 20
$ws.Range("F38").Value = 1215
$ws.Range("F39").Value = 1208
$ws.Range("F40").Value = 65
$ws.Range("F41").Value = 87
$ws.Range("F42").Value = 223
$ws.Range("F43").Value = 1490
$ws.Range("F44").Value = 2393
$ws.Range("F46").Value = 902
$ws.Range("F47").Value = 282
$ws.Range("F48").Value = 1244
$ws.Range("F49").Value = 23
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 676
$ws.Range("F3").Value = 873
$ws.Range("F4").Value = 94
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 676
$ws.Range("F5").Value = 873
$ws.Range("F6").Value = 94
$ws.Range("F9").Value = 154
$ws.Range("F11").Value = 156
$ws.Range("F12").Value = 754
$ws.Range("F13").Value = 4125
$ws.Range("F14").Value = 4125
$ws.Range("F16").Value = 166
$ws.Range("F18").Value = 5892
$ws.Range("F20").Value = 2284
$ws.Range("F22").Value = 157
$ws.Range("F24").Value = 8862
$ws.Range("F27").Value = 1696
$ws.Range("F28").Value = 2277
$ws.Range("F29").Value = 2371
$ws.Range("F30").Value = 1373
$ws.Range("F32").Value = 1922
$ws.Range("F34").Value = 52
$ws.Range("F37").Value = 37
$ws.Range("F39").Value = 1215
$ws.Range("F40").Value = 65
$ws.Range("F41").Value = 223
$ws.Range("F42").Value = 1490
$ws.Range("F43").Value = 2393
$ws.Range("F44").Value = 902
$ws.Range("F46").Value = 282
$ws.Range("F50").Value = 1244
